# Apply "Add data for 2022-05-20" (carjacking-by-neighborhood-by-month.xlsx)
#  - Rename sheet / update header & title text: "through May 09" -> "through May 12"
#  - Update / add per-neighborhood counts for the "May 2022" column (B) and a handful
#    of other month columns for specific neighborhoods (rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab and update the column-header / title cell text ---
$ws.Name = "Through 2022-05-12"
$ws.Range("B1").Value = "May 2022 (through May 12)"

# --- Cell value updates (both modifications of existing values and additions of
#     previously-empty cells) ---
$updates = @{
    "G2"  = 2
    "AK2" = 1
    "B3"  = 2
    "G3"  = 4
    "L3"  = 2
    "Q3"  = 3
    "AF3" = 2
    "V4"  = 1
    "AF4" = 2
    "G5"  = 1
    "V5"  = 4
    "Q7"  = 2
    "L10" = 1
    "B11" = 4
    "L11" = 1
    "Q11" = 2
    "AA11" = 1
    "AA13" = 2
    "AF13" = 1
    "G15" = 1
    "AA21" = 3
    "G23" = 2
    "V25" = 2
    "L26" = 1
    "AA26" = 1
    "V28" = 1
    "G29" = 1
    "B30" = 2
    "AA30" = 1
    "AA34" = 1
    "V35" = 1
    "AA37" = 1
    "B38" = 1
    "G39" = 3
    "B45" = 1
    "G46" = 3
    "V46" = 1
    "AF48" = 1
    "G63" = 2
    "L71" = 1
    "B73" = 1
    "B74" = 1
    "L89" = 1
    "B93" = 4
    "G93" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
